$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear cell A4 which held "Jorge Riopedre Vega"
$ws.Range("A4").ClearContents()

# Update the selection to A4 instead of C4
$ws.Range("A4").Select()
